# Apply changes: update header text and append cents to Price Per Unit values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for column D
$ws.Range("D1").Value = "Price Per Unit (£s)"

# Update the price values (row -> new value)
$priceUpdates = @{
    2  = 1054.63
    3  = 1629.57
    4  = 1606.67
    5  = 1201.03
    6  = 1434.28
    7  = 289.95
    8  = 500.08
    9  = 2100.87
    10 = 132.42
    11 = 994.4299999999999
    12 = 896.17
    13 = 905.11
    14 = 1076.55
    15 = 1592.19
    16 = 325.24
    17 = 1058.77
    18 = 291.21
    19 = 1255.21
    20 = 2135.78
    21 = 497.25
    22 = 1129.05
    23 = 854.6900000000001
    24 = 830.78
    25 = 1131.27
    26 = 2449.87
    27 = 1578.02
    28 = 2297.42
    29 = 2606.77
    30 = 277.81
    31 = 2172.67
    32 = 2072.84
    33 = 888.53
    34 = 618.09
    35 = 644.99
    36 = 2075.84
    37 = 2581.93
    38 = 724.99
    39 = 276.75
    40 = 1432.77
    41 = 734.76
    42 = 39.34
    43 = 2900.69
    44 = 1199.25
    45 = 847.83
    46 = 650.3099999999999
    47 = 1495.55
    48 = 2106.64
    49 = 2213.31
    50 = 1038.88
    51 = 1595.77
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Cells.Item($row, 4).Value = $priceUpdates[$row]
}
